$wb = $excel.ActiveWorkbook

# --- Update shared string content (prefix table: "url" header -> "uri") ---
$wsPrefixes = $wb.Worksheets.Item("_prefixes")
$wsPrefixes.Range("B1").Value = "uri"

# --- Update custom vocabulary table: "foafAgent" -> "foaf:Agent" ---
$wsCustomVoc = $wb.Worksheets.Item("_customVoc")
$wsCustomVoc.Range("B2").Value = "foaf:Agent"

# --- Update selection on _prefixes sheet (no longer the active tab) ---
$wsPrefixes.Range("E8").Select() | Out-Null

# --- Make _customVoc the active sheet/tab, with a new selection ---
$wsCustomVoc.Activate() | Out-Null
$wsCustomVoc.Range("B3").Select() | Out-Null
